$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - Dengue
$ws.Range("C6").Value = 3

# Row 7 - Defectos congenitos
$ws.Range("C7").Value = 3
$ws.Range("E7").Value = 0.22

# Row 11 - Agresiones por animales potencialmente transmisores de rabia
$ws.Range("C11").Value = 41

# Row 12 - Hepatitis a
$ws.Range("C12").Value = 1
$ws.Range("E12").Value = 0.37

# Row 14 - Enfermedades huerfanas - raras
$ws.Range("C14").Value = 2
$ws.Range("E14").Value = 0.04

# Row 18 - Intento de suicidio
$ws.Range("C18").Value = 9

# Row 19 - Iad - infecciones asociadas a dispositivos - individual
$ws.Range("C19").Value = 1
$ws.Range("E19").Value = 0.37

# Row 20 - Intoxicaciones
$ws.Range("C20").Value = 7
$ws.Range("E20").Value = 0.09

# Row 29 - Parotiditis
$ws.Range("C29").Value = 0
$ws.Range("E29").Value = 1

# Row 30 - Sifilis congenita
$ws.Range("C30").Value = 0
$ws.Range("E30").Value = 1

# Row 31 - Sifilis gestacional
$ws.Range("C31").Value = 1
$ws.Range("E31").Value = 0.37

# Row 32 - Tuberculosis
$ws.Range("C32").Value = 6
$ws.Range("E32").Value = 0.04

# Row 33 - Varicela individual
$ws.Range("C33").Value = 4
$ws.Range("E33").Value = 0.15

# Row 34 - Vih/sida/mortalidad por sida
$ws.Range("C34").Value = 8
$ws.Range("E34").Value = 0.14
